# Selenium package update and bug fixes:
# Supplier name changed in release -> set the "Supplier" column (N) value to
# "Mystifly" on every data row of all three sheets, and update the
# selection/scroll state to match the recorded view at commit time.

$wb = $excel.ActiveWorkbook

# --- Air_Mystifly_OneWay ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("Air_Mystifly_OneWay")
$ws1.Range("N2").Value = "Mystifly"
$ws1.Range("N3").Value = "Mystifly"
$ws1.Range("N4").Value = "Mystifly"
$ws1.Range("N5").Value = "Mystifly"
$ws1.Range("N5").Select()

# --- Air_Mystifly_RoundTrip -------------------------------------------------
$ws2 = $wb.Worksheets.Item("Air_Mystifly_RoundTrip")
$ws2.Range("N2").Value = "Mystifly"
$ws2.Range("N2").Select()

# --- Air_Mystifly_Multicity -------------------------------------------------
$ws3 = $wb.Worksheets.Item("Air_Mystifly_Multicity")
$ws3.Range("N2").Value = "Mystifly"
$ws3.Range("N2").Select()
